$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new descriptive rows (MCH175-1 / MCH175-2) beneath the header ---
$ws.Range("A2").Value = "MCH175-1"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 23F | GRAP COUNT NUMER: NONE"

$ws.Range("A3").Value = "MCH175-2"
$ws.Range("E3").Value = "Series"
$ws.Range("F3").Value = "1 Box"
$ws.Range("G3").Value = "LOCATION: 23F | GRAP COUNT NUMER: NONE"

# --- Apply the "Calibri 10" body font (theme text color) used throughout the new rows ---
$ws.Range("A2").Font.ThemeColor = 1
$ws.Range("A2").Font.Name = "Calibri"

# Propagate that same cell format (copy/paste) to the rest of the plain cells,
# including the blank-but-present ones, without touching their values.
$ws.Range("A2").Copy()
$ws.Range("C2:E2").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("G2:H2").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("C3:E3").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("G3:H3").PasteSpecial(-4122)

# Re-apply the literal values (PasteSpecial(formats) above must not disturb them,
# but make sure the text survives identically).
$ws.Range("A2").Value = "MCH175-1"
$ws.Range("E2").Value = "Series"
$ws.Range("G2").Value = "LOCATION: 23F | GRAP COUNT NUMER: NONE"
$ws.Range("A3").Value = "MCH175-2"
$ws.Range("E3").Value = "Series"
$ws.Range("G3").Value = "LOCATION: 23F | GRAP COUNT NUMER: NONE"

# --- The extent/medium column (F) uses a slightly different cell format (same font) ---
$ws.Range("F2").Font.ThemeColor = 1
$ws.Range("F2").Font.Name = "Calibri"
$ws.Range("F2").Value = "1 Box"
$ws.Range("F2").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("F3").Value = "1 Box"

$excel.CutCopyMode = $false

# --- Restore the frozen header pane / selection over the new data range ---
$ws.Range("A2:I3").Select()
$excel.ActiveWindow.FreezePanes = $true
